{"js": "// The site rebuild dropped the trailing footer block: a blank paragraph,\n// a blank page-break paragraph, and the \"\u00a9 2020 ... Creative Commons\n// Attribution\" paragraph that followed the bibliography text. Find that\n// copyright paragraph by its text, confirm the two paragraphs in front of\n// it are the expected blank ones, and delete all three so the bibliography\n// paragraph is followed directly by the document's final (blank,\n// page-break) paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"\u00a9 2020\";\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  const copyrightPara = paragraphs.items[targetIndex];\n  const blank1 = targetIndex - 1 >= 0 ? paragraphs.items[targetIndex - 1] : null; // page-break blank\n  const blank2 = targetIndex - 2 >= 0 ? paragraphs.items[targetIndex - 2] : null; // plain blank\n\n  if (blank1) blank1.paragraphFormat.load(\"pageBreakBefore\");\n  await context.sync();\n\n  const toDelete = [copyrightPara];\n  // Only remove the preceding paragraphs when they are genuinely blank,\n  // matching the shape of the removed block (blank, blank+page-break,\n  // copyright text) rather than assuming a fixed layout.\n  if (blank1 && blank1.text === \"\" && blank1.paragraphFormat.pageBreakBefore) {\n    toDelete.push(blank1);\n    if (blank2 && blank2.text === \"\") {\n      toDelete.push(blank2);\n    }\n  }\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# The site rebuild dropped the trailing footer block: a blank paragraph,\n# a blank page-break paragraph, and the \"\u00a9 2020 ... Creative Commons\n# Attribution\" paragraph that followed the bibliography text. Locate that\n# copyright paragraph by its (ASCII-safe) text, confirm the two paragraphs\n# immediately before it are the expected blank ones, and delete all three -\n# leaving the bibliography paragraph followed directly by the document's\n# final (blank, page-break) paragraph.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"Creative Commons Attribution\")) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ge 1) {\n    $deleteIndexes = @($targetIndex)\n\n    if ($targetIndex - 1 -ge 1) {\n        $blank1 = $d.Paragraphs.Item($targetIndex - 1)\n        $blank1IsEmpty = $blank1.Range.Text.Length -le 1\n        $blank1HasPageBreak = $blank1.Format.PageBreakBefore\n\n        if ($blank1IsEmpty -and $blank1HasPageBreak) {\n            $deleteIndexes += ($targetIndex - 1)\n\n            if ($targetIndex - 2 -ge 1) {\n                $blank2 = $d.Paragraphs.Item($targetIndex - 2)\n                if ($blank2.Range.Text.Length -le 1) {\n                    $deleteIndexes += ($targetIndex - 2)\n                }\n            }\n        }\n    }\n\n    # Delete from the highest index down so earlier indexes stay valid.\n    $deleteIndexes = $deleteIndexes | Sort-Object -Descending\n    foreach ($idx in $deleteIndexes) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
